$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-24 22:18:40'
$ws.Range('H2').NumberFormat = '@'
$ws.Range('H2').Value = '42%'
$ws.Range('E3').Value = '2026-02-24 22:18:43'
$ws.Range('E4').Value = '2026-02-24 22:18:46'
$ws.Range('H4').NumberFormat = '@'
$ws.Range('H4').Value = '73%'
$ws.Range('O4').Value = '12.6 °C'
$ws.Range('E5').Value = '2026-02-24 22:18:48'
$ws.Range('O5').Value = '5.8 °C'
$ws.Range('E6').Value = '2026-02-24 22:18:51'
$ws.Range('H6').NumberFormat = '@'
$ws.Range('H6').Value = '71%'
$ws.Range('O6').Value = '13.9 °C'
$ws.Range('E7').Value = '2026-02-24 22:18:54'
$ws.Range('H7').NumberFormat = '@'
$ws.Range('H7').Value = '73%'
$ws.Range('O7').Value = '14.0 °C'
$ws.Range('E8').Value = '2026-02-24 22:18:57'
$ws.Range('H8').NumberFormat = '@'
$ws.Range('H8').Value = '45%'
$ws.Range('E9').Value = '2026-02-24 22:18:59'
$ws.Range('O9').Value = '11.7 °C'
$ws.Range('E10').Value = '2026-02-24 22:19:02'
$ws.Range('H10').NumberFormat = '@'
$ws.Range('H10').Value = '78%'
$ws.Range('O10').Value = '10.9 °C'
$ws.Range('E11').Value = '2026-02-24 22:19:04'
$ws.Range('O11').Value = '9.0 °C'
$ws.Range('E12').Value = '2026-02-24 22:19:07'
$ws.Range('O12').Value = '10.4 °C'
$ws.Range('E13').Value = '2026-02-24 22:19:10'
$ws.Range('K13').Value = '14.8 MJ/m2'
$ws.Range('E14').Value = '2026-02-24 22:19:12'
$ws.Range('N14').Value = '5.7 °C 21:59 TU'
$ws.Range('O14').Value = '11.3 °C'
$ws.Range('E15').Value = '2026-02-24 22:19:15'
$ws.Range('O15').Value = '11.7 °C'
$ws.Range('E16').Value = '2026-02-24 22:19:17'
$ws.Range('E17').Value = '2026-02-24 22:19:20'
$ws.Range('H17').NumberFormat = '@'
$ws.Range('H17').Value = '32%'
$ws.Range('E18').Value = '2026-02-24 22:19:23'
$ws.Range('H18').NumberFormat = '@'
$ws.Range('H18').Value = '77%'
$ws.Range('O18').Value = '11.1 °C'
$ws.Range('E19').Value = '2026-02-24 22:19:25'
$ws.Range('O19').Value = '12.5 °C'
$ws.Range('E20').Value = '2026-02-24 22:19:28'
$ws.Range('O20').Value = '3.7 °C'
$ws.Range('E21').Value = '2026-02-24 22:19:31'
$ws.Range('J21').Value = '1021.9 hPa'
$ws.Range('O21').Value = '9.7 °C'
$ws.Range('E22').Value = '2026-02-24 22:19:34'
$ws.Range('L22').Value = '20.9 km/h - 308º 21:54 TU'
$ws.Range('O22').Value = '3.4 °C'
$ws.Range('E23').Value = '2026-02-24 22:19:36'
$ws.Range('N23').Value = '2.9 °C 21:59 TU'
$ws.Range('E24').Value = '2026-02-24 22:19:39'
$ws.Range('J24').Value = '1021.0 hPa'
$ws.Range('L24').Value = '14.8 km/h - 73º 21:33 TU'
$ws.Range('E25').Value = '2026-02-24 22:19:42'
$ws.Range('E26').Value = '2026-02-24 22:19:44'
$ws.Range('J26').Value = '1018.8 hPa'
$ws.Range('O26').Value = '11.6 °C'
$ws.Range('E27').Value = '2026-02-24 22:19:47'
$ws.Range('E28').Value = '2026-02-24 22:19:50'
$ws.Range('O28').Value = '11.4 °C'
$ws.Range('E29').Value = '2026-02-24 22:19:52'
$ws.Range('H29').NumberFormat = '@'
$ws.Range('H29').Value = '90%'
$ws.Range('E30').Value = '2026-02-24 22:19:55'
$ws.Range('J30').Value = '1019.6 hPa'
$ws.Range('E31').Value = '2026-02-24 22:19:58'
$ws.Range('N31').Value = '13.1 °C 21:59 TU'
$ws.Range('O31').Value = '15.7 °C'
$ws.Range('E32').Value = '2026-02-24 22:20:01'
$ws.Range('O32').Value = '6.9 °C'
$ws.Range('E33').Value = '2026-02-24 22:20:03'
$ws.Range('O33').Value = '8.5 °C'
$ws.Range('E34').Value = '2026-02-24 22:20:06'
$ws.Range('O34').Value = '4.6 °C'
$ws.Range('E35').Value = '2026-02-24 22:20:09'
$ws.Range('E36').Value = '2026-02-24 22:20:11'
$ws.Range('O36').Value = '12.8 °C'
$ws.Range('E37').Value = '2026-02-24 22:20:14'
$ws.Range('O37').Value = '8.5 °C'
$ws.Range('E38').Value = '2026-02-24 22:20:17'
$ws.Range('O38').Value = '11.8 °C'
$ws.Range('E39').Value = '2026-02-24 22:20:19'
$ws.Range('H39').NumberFormat = '@'
$ws.Range('H39').Value = '36%'
$ws.Range('L39').Value = '23.4 km/h - 251º 21:54 TU'
$ws.Range('E40').Value = '2026-02-24 22:20:22'
$ws.Range('O40').Value = '8.4 °C'
$ws.Range('E41').Value = '2026-02-24 22:20:24'
$ws.Range('J41').Value = '1020.4 hPa'
$ws.Range('E42').Value = '2026-02-24 22:20:27'
$ws.Range('E43').Value = '2026-02-24 22:20:30'
$ws.Range('E44').Value = '2026-02-24 22:20:32'
$ws.Range('E45').Value = '2026-02-24 22:20:35'
$ws.Range('E46').Value = '2026-02-24 22:20:38'
$ws.Range('J46').Value = '1021.0 hPa'
